# Insert a new record at row 143 (pushes the existing block of rows
# 143-258 down to 144-259, preserving all of their values/styles), then
# populate the newly inserted row with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(143).Insert()

$ws.Cells.Item(143, 1).Value = 5
$ws.Cells.Item(143, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(143, 3).Value = "Maule"
$ws.Cells.Item(143, 4).Value = 44741
$ws.Cells.Item(143, 5).Value = 7
$ws.Cells.Item(143, 6).Value = 100112008
$ws.Cells.Item(143, 7).Value = "Coliflor"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 4000
$ws.Cells.Item(143, 11).Value = 900
$ws.Cells.Item(143, 12).Value = 900
$ws.Cells.Item(143, 13).Value = 900
$ws.Cells.Item(143, 14).Value = "$/unidad"
$ws.Cells.Item(143, 15).Value = "Región del Maule"
$ws.Cells.Item(143, 16).Value = 900
$ws.Cells.Item(143, 17).Value = 1
$ws.Cells.Item(143, 18).Value = "Hortaliza"
